# Predicting car purchase decision.docx - apply commit changes
$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $ok = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        Write-Output ("WARNING: replace failed for: " + $findText)
    }
    return $ok
}

# 1) "Define independent (X) and dependent (y) variables" -> "Defining Independent Variables"
Replace-Text "Define independent (X) and dependent (y) variables" "Defining Independent Variables" | Out-Null

# 2) LDA predictors paragraph: insert "Gender, " before "Age and Annual Salary"
Replace-Text "The LDA model uses Age and Annual Salary as the independent variables" "The LDA model uses Gender, Age and Annual Salary as the independent variables" | Out-Null

# 3) Confusion matrix numbers: 38 -> 41, 90 -> 87
Replace-Text " [ 38  90]]" " [ 41  87]]" | Out-Null
Replace-Text "False Negatives (FN): 38 " "False Negatives (FN): 41 " | Out-Null
Replace-Text "True Positives (TP): 90 " "True Positives (TP): 87 " | Out-Null

# 4) Predicted counts: 104 -> 101 ; 196 -> 199
Replace-Text "Predicted Purchased: 104 users" "Predicted Purchased: 101 users" | Out-Null
Replace-Text "Predicted Not Purchased: 196 users" "Predicted Not Purchased: 199 users" | Out-Null

# 5) Insights false negatives parenthetical: 38 -> 41
Replace-Text "While the number of false negatives (38) suggests" "While the number of false negatives (41) suggests" | Out-Null

# 6) Insights predictors sentence
Replace-Text "The use of only two predictors (Age and Salary) achieved" "The use of predictors (Gender, Age and Salary) achieved" | Out-Null

Write-Output "Done with text replacements"
Write-Output $d.Content.Text
